$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v_A2 = @'
<Policy PolicyId="mli-interface" RuleCombiningAlgId="urn:oasis:names:tc:xacml:3.0:rule-combining-algorithm:permit-overrides" Version="1.0">
        <Target>
            <AnyOf>
                <AllOf>
                    <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                        <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">VIO</AttributeValue>
                        <AttributeDesignator AttributeId="http://authz-interop.org/AAA/xacml/subject/subject-role" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                    </Match>
                </AllOf>
            </AnyOf>
        </Target>
        <Rule Effect="Permit" RuleId="mli:vi-operations">
            <Target>
                <AnyOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">VI</AttributeValue>
                            <AttributeDesignator AttributeId="http://authz-interop.org/AAA/xacml/resource/resource-type" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                </AnyOf>
                <AnyOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:Request-VI</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:Instantiate-VI</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:Decommission-VI</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:ReplanningVI:Modify-Time</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                </AnyOf>
            </Target>
        </Rule>
        <Rule Effect="Permit" RuleId="mli:vi-request-operations">
            <Target>
                <AnyOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">VI-Request</AttributeValue>
                            <AttributeDesignator AttributeId="http://authz-interop.org/AAA/xacml/resource/resource-type" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                </AnyOf>
                <AnyOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:Query-VI-Request-Status</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:Get-SLA-Offer</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:Sign-SLA-Offer</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                </AnyOf>
            </Target>
        </Rule>
        <Rule Effect="Permit" RuleId="mli:replanning:vr-it-operations">
            <Target>
                <AnyOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">VR-IT</AttributeValue>
                            <AttributeDesignator AttributeId="http://authz-interop.org/AAA/xacml/resource/resource-type" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                </AnyOf>
                <AnyOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:ReplanningVI:Add-VR-IT</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:ReplanningVI:Modify-VR-IT</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:ReplanningVI:Delete-VR-IT</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                </AnyOf>
            </Target>
        </Rule>
        <Rule Effect="Permit" RuleId="mli:replanning:vlink-operations">
            <Target>
                <AnyOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">VLink</AttributeValue>
                            <AttributeDesignator AttributeId="http://authz-interop.org/AAA/xacml/resource/resource-type" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                </AnyOf>
                <AnyOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:ReplanningVI:Add-VLink</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:ReplanningVI:Modify-VLink</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                    <AllOf>
                        <Match MatchId="urn:oasis:names:tc:xacml:1.0:function:string-equal">
                            <AttributeValue DataType="http://www.w3.org/2001/XMLSchema#string">MLI:ReplanningVI:Delete-VLink</AttributeValue>
                            <AttributeDesignator AttributeId="urn:oasis:names:tc:xacml:1.0:action:action-id" DataType="http://www.w3.org/2001/XMLSchema#string" MustBePresent="false"/>
                        </Match>
                    </AllOf>
                </AnyOf>
            </Target>
        </Rule>
    </Policy>
'@
$ws.Range("A2").Value = $v_A2

$v_C2 = @'
{
    "datalog_subjects": "Subject(S), has_role(S, 'VIO') :- Subject(S).",
    "datalog_objects": "Resource(R), VI(R) :- Resource(R), has_type(R, 'VI'). Resource(R), VI_Request(R) :- Resource(R), has_type(R, 'VI-Request'). Resource(R), VR_IT(R) :- Resource(R), has_type(R, 'VR-IT'). Resource(R), VLink(R) :- Resource(R), has_type(R, 'VLink').",
    "datalog_relationships": "has_role(S, 'VIO') :- Subject(S). has_type(R, 'VI') :- VI(R). has_type(R, 'VI-Request') :- VI_Request(R). has_type(R, 'VR-IT') :- VR_IT(R). has_type(R, 'VLink') :- VLink(R).",
    "datalog_actions": "can_perform(S, 'MLI:Request-VI', R) :- Subject(S), VI(R), has_role(S, 'VIO'). can_perform(S, 'MLI:Instantiate-VI', R) :- Subject(S), VI(R), has_role(S, 'VIO'). can_perform(S, 'MLI:Decommission-VI', R) :- Subject(S), VI(R), has_role(S, 'VIO'). can_perform(S, 'MLI:ReplanningVI:Modify-Time', R) :- Subject(S), VI(R), has_role(S, 'VIO'). can_perform(S, 'MLI:Query-VI-Request-Status', R) :- Subject(S), VI_Request(R), has_role(S, 'VIO'). can_perform(S, 'MLI:Get-SLA-Offer', R) :- Subject(S), VI_Request(R), has_role(S, 'VIO'). can_perform(S, 'MLI:Sign-SLA-Offer', R) :- Subject(S), VI_Request(R), has_role(S, 'VIO'). can_perform(S, 'MLI:ReplanningVI:Add-VR-IT', R) :- Subject(S), VR_IT(R), has_role(S, 'VIO'). can_perform(S, 'MLI:ReplanningVI:Modify-VR-IT', R) :- Subject(S), VR_IT(R), has_role(S, 'VIO'). can_perform(S, 'MLI:ReplanningVI:Delete-VR-IT', R) :- Subject(S), VR_IT(R), has_role(S, 'VIO'). can_perform(S, 'MLI:ReplanningVI:Add-VLink', R) :- Subject(S), VLink(R), has_role(S, 'VIO'). can_perform(S, 'MLI:ReplanningVI:Modify-VLink', R) :- Subject(S), VLink(R), has_role(S, 'VIO'). can_perform(S, 'MLI:ReplanningVI:Delete-VLink', R) :- Subject(S), VLink(R), has_role(S, 'VIO')."
}
'@
$ws.Range("C2").Value = $v_C2

$v_C3 = @'
{
    "datalog_subjects": "Subject(S), Role(S, 'VIO') :- Subject(S).",
    "datalog_objects": "Resource(R), ResourceType(R, 'VR') :- Resource(R). ResourceType(R, 'VR-RP') :- Resource(R). ResourceType(R, 'VR-RP-Info') :- Resource(R). ResourceType(R, 'VR-State-Info') :- Resource(R). ResourceType(R, 'VR-Power-Info') :- Resource(R). ResourceType(R, 'VR-Status-Info') :- Resource(R). ResourceType(R, 'VR-Info') :- Resource(R).",
    "datalog_relationships": "",
    "datalog_actions": "can_perform(S, R, 'SLI:Decommission-VR-IT') :- Subject(S), Resource(R), Role(S, 'VIO'), ResourceType(R, 'VR'). can_perform(S, R, 'SLI:Operate-VR') :- Subject(S), Resource(R), Role(S, 'VIO'), ResourceType(R, 'VR'). can_perform(S, R, 'SLI:Add-VirtualNetworkIf') :- Subject(S), Resource(R), Role(S, 'VIO'), ResourceType(R, 'VR'). can_perform(S, R, 'SLI:Remove-VirtualNetworkIf') :- Subject(S), Resource(R), Role(S, 'VIO'), ResourceType(R, 'VR'). can_perform(S, R, 'SLI:Create-StorageImage') :- Subject(S), Resource(R), Role(S, 'VIO'), ResourceType(R, 'VR'). can_perform(S, R, 'SLI:Remove-StorageImage') :- Subject(S), Resource(R), Role(S, 'VIO'), ResourceType(R, 'VR'). can_perform(S, R, 'SLI:Instantiate-VR-IT') :- Subject(S), Resource(R), Role(S, 'VIO'), ResourceType(R, 'VR-RP'). can_perform(S, R, 'SLI:Get-Available-VR-Pool-IT') :- Subject(S), Resource(R), Role(S, 'VIO'), ResourceType(R, 'VR-RP-Info'). can_perform(S, R, 'SLI:Monitor-VR-Info') :- Subject(S), Resource(R), Role(S, 'VIO'), ResourceType(R, 'VR-State-Info'), ResourceType(R, 'VR-Power-Info'), ResourceType(R, 'VR-Status-Info'). can_perform(S, R, 'SLI:Subscribe-VR-Monitoring') :- Subject(S), Resource(R), Role(S, 'VIO'), ResourceType(R, 'VR-Info'). can_perform(S, R, 'SLI:Unsubscribe-VR-Monitoring') :- Subject(S), Resource(R), Role(S, 'VIO'), ResourceType(R, 'VR-Info')."
}
'@
$ws.Range("C3").Value = $v_C3

$v_C4 = @'
{
    "datalog_subjects": "Subject(S).",
    "datalog_objects": "Resource(R).",
    "datalog_relationships": "has_role(S, 'VIO') :- Subject(S).",
    "datalog_actions": "can_access(S, R, A) :- Subject(S), Resource(R), has_role(S, 'VIO'), resource_type(R, 'VI'), action_in(A, ['MLI:Request-VI', 'MLI:Instantiate-VI', 'MLI:Decommission-VI', 'MLI:ReplanningVI:Modify-Time']).\ncan_access(S, R, A) :- Subject(S), Resource(R), has_role(S, 'VIO'), resource_type(R, 'VI-Request'), action_in(A, ['MLI:Query-VI-Request-Status', 'MLI:Get-SLA-Offer', 'MLI:Sign-SLA-Offer']).\ncan_access(S, R, A) :- Subject(S), Resource(R), has_role(S, 'VIO'), resource_type(R, 'VR-IT'), action_in(A, ['MLI:ReplanningVI:Add-VR-IT', 'MLI:ReplanningVI:Modify-VR-IT', 'MLI:ReplanningVI:Delete-VR-IT']).\ncan_access(S, R, A) :- Subject(S), Resource(R), has_role(S, 'VIO'), resource_type(R, 'VLink'), action_in(A, ['MLI:ReplanningVI:Add-VLink', 'MLI:ReplanningVI:Modify-VLink', 'MLI:ReplanningVI:Delete-VLink'])."
}
'@
$ws.Range("C4").Value = $v_C4

$v_C5 = @'
{
    "datalog_subjects": "VIO(S).",
    "datalog_objects": "VR(R), VR_RP(RP), VR_RP_Info(RI), VR_State_Info(SI), VR_Power_Info(PI), VR_Status_Info(STI), VR_Info(VI).",
    "datalog_relationships": "",
    "datalog_actions": "can_perform(S, R, A) :- VIO(S), VR(R), action_in(A, [\"SLI:Decommission-VR-IT\", \"SLI:Operate-VR\", \"SLI:Add-VirtualNetworkIf\", \"SLI:Remove-VirtualNetworkIf\", \"SLI:Create-StorageImage\", \"SLI:Remove-StorageImage\"]). can_perform(S, RP, \"SLI:Instantiate-VR-IT\") :- VIO(S), VR_RP(RP). can_perform(S, RI, \"SLI:Get-Available-VR-Pool-IT\") :- VIO(S), VR_RP_Info(RI). can_perform(S, SI, \"SLI:Monitor-VR-Info\") :- VIO(S), VR_State_Info(SI). can_perform(S, PI, \"SLI:Monitor-VR-Info\") :- VIO(S), VR_Power_Info(PI). can_perform(S, STI, \"SLI:Monitor-VR-Info\") :- VIO(S), VR_Status_Info(STI). can_perform(S, VI, \"SLI:Subscribe-VR-Monitoring\") :- VIO(S), VR_Info(VI). can_perform(S, VI, \"SLI:Unsubscribe-VR-Monitoring\") :- VIO(S), VR_Info(VI)."
}
'@
$ws.Range("C5").Value = $v_C5

$v_C6 = @'
{
    "datalog_subjects": "Subject(S), has_role(S, 'VIO') :- Subject(S).",
    "datalog_objects": "Resource(R), has_type(R, 'VI') :- Resource(R). has_type(R, 'VI-Request') :- Resource(R). has_type(R, 'VR-IT') :- Resource(R). has_type(R, 'VLink') :- Resource(R).",
    "datalog_relationships": "",
    "datalog_actions": "can_perform(S, R, 'MLI:Request-VI') :- Subject(S), Resource(R), has_role(S, 'VIO'), has_type(R, 'VI'). can_perform(S, R, 'MLI:Instantiate-VI') :- Subject(S), Resource(R), has_role(S, 'VIO'), has_type(R, 'VI'). can_perform(S, R, 'MLI:Decommission-VI') :- Subject(S), Resource(R), has_role(S, 'VIO'), has_type(R, 'VI'). can_perform(S, R, 'MLI:ReplanningVI:Modify-Time') :- Subject(S), Resource(R), has_role(S, 'VIO'), has_type(R, 'VI'). can_perform(S, R, 'MLI:Query-VI-Request-Status') :- Subject(S), Resource(R), has_role(S, 'VIO'), has_type(R, 'VI-Request'). can_perform(S, R, 'MLI:Get-SLA-Offer') :- Subject(S), Resource(R), has_role(S, 'VIO'), has_type(R, 'VI-Request'). can_perform(S, R, 'MLI:Sign-SLA-Offer') :- Subject(S), Resource(R), has_role(S, 'VIO'), has_type(R, 'VI-Request'). can_perform(S, R, 'MLI:ReplanningVI:Add-VR-IT') :- Subject(S), Resource(R), has_role(S, 'VIO'), has_type(R, 'VR-IT'). can_perform(S, R, 'MLI:ReplanningVI:Modify-VR-IT') :- Subject(S), Resource(R), has_role(S, 'VIO'), has_type(R, 'VR-IT'). can_perform(S, R, 'MLI:ReplanningVI:Delete-VR-IT') :- Subject(S), Resource(R), has_role(S, 'VIO'), has_type(R, 'VR-IT'). can_perform(S, R, 'MLI:ReplanningVI:Add-VLink') :- Subject(S), Resource(R), has_role(S, 'VIO'), has_type(R, 'VLink'). can_perform(S, R, 'MLI:ReplanningVI:Modify-VLink') :- Subject(S), Resource(R), has_role(S, 'VIO'), has_type(R, 'VLink'). can_perform(S, R, 'MLI:ReplanningVI:Delete-VLink') :- Subject(S), Resource(R), has_role(S, 'VIO'), has_type(R, 'VLink')."
}
'@
$ws.Range("C6").Value = $v_C6

$v_C7 = @'
{
    "datalog_subjects": "Subject(S), VIO(S) :- Subject(S), subject_role(S, 'VIO').",
    "datalog_objects": "Resource(R), VR(R) :- Resource(R), resource_type(R, 'VR').\nResource(R), VR_RP(R) :- Resource(R), resource_type(R, 'VR-RP').\nResource(R), VR_RP_Info(R) :- Resource(R), resource_type(R, 'VR-RP-Info').\nResource(R), VR_State_Info(R) :- Resource(R), resource_type(R, 'VR-State-Info').\nResource(R), VR_Power_Info(R) :- Resource(R), resource_type(R, 'VR-Power-Info').\nResource(R), VR_Status_Info(R) :- Resource(R), resource_type(R, 'VR-Status-Info').\nResource(R), VR_Info(R) :- Resource(R), resource_type(R, 'VR-Info').",
    "datalog_relationships": "",
    "datalog_actions": "can_perform(S, R, 'SLI:Decommission-VR-IT') :- VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Operate-VR') :- VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Add-VirtualNetworkIf') :- VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Remove-VirtualNetworkIf') :- VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Create-StorageImage') :- VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Remove-StorageImage') :- VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Instantiate-VR-IT') :- VIO(S), VR_RP(R).\ncan_perform(S, R, 'SLI:Get-Available-VR-Pool-IT') :- VIO(S), VR_RP_Info(R).\ncan_perform(S, R, 'SLI:Monitor-VR-Info') :- VIO(S), VR_State_Info(R), VR_Power_Info(R), VR_Status_Info(R).\ncan_perform(S, R, 'SLI:Subscribe-VR-Monitoring') :- VIO(S), VR_Info(R).\ncan_perform(S, R, 'SLI:Unsubscribe-VR-Monitoring') :- VIO(S), VR_Info(R)."
}
'@
$ws.Range("C7").Value = $v_C7

$v_C8 = @'
{
    "datalog_subjects": "Subject(S), has_role(S, VIO) :- Subject(S).",
    "datalog_objects": "Resource(R), has_type(R, VI) :- Resource(R). has_type(R, VI_Request) :- Resource(R). has_type(R, VR_IT) :- Resource(R). has_type(R, VLink) :- Resource(R).",
    "datalog_relationships": "",
    "datalog_actions": "can_access(S, R, Action) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VI), Action = 'MLI:Request-VI'; Action = 'MLI:Instantiate-VI'; Action = 'MLI:Decommission-VI'; Action = 'MLI:ReplanningVI:Modify-Time'. can_access(S, R, Action) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VI_Request), Action = 'MLI:Query-VI-Request-Status'; Action = 'MLI:Get-SLA-Offer'; Action = 'MLI:Sign-SLA-Offer'. can_access(S, R, Action) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VR_IT), Action = 'MLI:ReplanningVI:Add-VR-IT'; Action = 'MLI:ReplanningVI:Modify-VR-IT'; Action = 'MLI:ReplanningVI:Delete-VR-IT'. can_access(S, R, Action) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VLink), Action = 'MLI:ReplanningVI:Add-VLink'; Action = 'MLI:ReplanningVI:Modify-VLink'; Action = 'MLI:ReplanningVI:Delete-VLink'."
}
'@
$ws.Range("C8").Value = $v_C8

$v_C9 = @'
{
    "datalog_subjects": "Subject(S), VIO(S).",
    "datalog_objects": "Resource(R), VR(R), VR_RP(R), VR_RP_Info(R), VR_State_Info(R), VR_Power_Info(R), VR_Status_Info(R), VR_Info(R).",
    "datalog_relationships": "",
    "datalog_actions": "can_perform(S, R, Action) :- Subject(S), VIO(S), Resource(R), Action(Action).\ncan_perform(S, R, 'SLI:Decommission-VR-IT') :- Subject(S), VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Operate-VR') :- Subject(S), VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Add-VirtualNetworkIf') :- Subject(S), VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Remove-VirtualNetworkIf') :- Subject(S), VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Create-StorageImage') :- Subject(S), VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Remove-StorageImage') :- Subject(S), VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Instantiate-VR-IT') :- Subject(S), VIO(S), VR_RP(R).\ncan_perform(S, R, 'SLI:Get-Available-VR-Pool-IT') :- Subject(S), VIO(S), VR_RP_Info(R).\ncan_perform(S, R, 'SLI:Monitor-VR-Info') :- Subject(S), VIO(S), VR_State_Info(R).\ncan_perform(S, R, 'SLI:Monitor-VR-Info') :- Subject(S), VIO(S), VR_Power_Info(R).\ncan_perform(S, R, 'SLI:Monitor-VR-Info') :- Subject(S), VIO(S), VR_Status_Info(R).\ncan_perform(S, R, 'SLI:Subscribe-VR-Monitoring') :- Subject(S), VIO(S), VR_Info(R).\ncan_perform(S, R, 'SLI:Unsubscribe-VR-Monitoring') :- Subject(S), VIO(S), VR_Info(R)."
}
'@
$ws.Range("C9").Value = $v_C9

$v_C10 = @'
{
    "datalog_subjects": "Subject(S), has_role(S, VIO) :- Subject(S).",
    "datalog_objects": "Resource(R), has_type(R, VI) :- Resource(R). has_type(R, VI_Request) :- Resource(R). has_type(R, VR_IT) :- Resource(R). has_type(R, VLink) :- Resource(R).",
    "datalog_relationships": "",
    "datalog_actions": "can_perform(S, R, MLI_Request_VI) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VI). can_perform(S, R, MLI_Instantiate_VI) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VI). can_perform(S, R, MLI_Decommission_VI) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VI). can_perform(S, R, MLI_ReplanningVI_Modify_Time) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VI). can_perform(S, R, MLI_Query_VI_Request_Status) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VI_Request). can_perform(S, R, MLI_Get_SLA_Offer) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VI_Request). can_perform(S, R, MLI_Sign_SLA_Offer) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VI_Request). can_perform(S, R, MLI_ReplanningVI_Add_VR_IT) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VR_IT). can_perform(S, R, MLI_ReplanningVI_Modify_VR_IT) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VR_IT). can_perform(S, R, MLI_ReplanningVI_Delete_VR_IT) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VR_IT). can_perform(S, R, MLI_ReplanningVI_Add_VLink) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VLink). can_perform(S, R, MLI_ReplanningVI_Modify_VLink) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VLink). can_perform(S, R, MLI_ReplanningVI_Delete_VLink) :- Subject(S), Resource(R), has_role(S, VIO), has_type(R, VLink)."
}
'@
$ws.Range("C10").Value = $v_C10

$v_C11 = @'
{
    "datalog_subjects": "VIO(S).",
    "datalog_objects": "VR(R), VR_RP(RP), VR_RP_Info(RI), VR_State_Info(SI), VR_Power_Info(PI), VR_Status_Info(TI), VR_Info(VI).",
    "datalog_relationships": "",
    "datalog_actions": "can_perform(S, R, 'SLI:Decommission-VR-IT') :- VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Operate-VR') :- VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Add-VirtualNetworkIf') :- VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Remove-VirtualNetworkIf') :- VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Create-StorageImage') :- VIO(S), VR(R).\ncan_perform(S, R, 'SLI:Remove-StorageImage') :- VIO(S), VR(R).\ncan_perform(S, RP, 'SLI:Instantiate-VR-IT') :- VIO(S), VR_RP(RP).\ncan_perform(S, RI, 'SLI:Get-Available-VR-Pool-IT') :- VIO(S), VR_RP_Info(RI).\ncan_perform(S, SI, 'SLI:Monitor-VR-Info') :- VIO(S), VR_State_Info(SI).\ncan_perform(S, PI, 'SLI:Monitor-VR-Info') :- VIO(S), VR_Power_Info(PI).\ncan_perform(S, TI, 'SLI:Monitor-VR-Info') :- VIO(S), VR_Status_Info(TI).\ncan_perform(S, VI, 'SLI:Subscribe-VR-Monitoring') :- VIO(S), VR_Info(VI).\ncan_perform(S, VI, 'SLI:Unsubscribe-VR-Monitoring') :- VIO(S), VR_Info(VI)."
}
'@
$ws.Range("C11").Value = $v_C11

$v_C12 = @'
{
    "datalog_subjects": "VIO(S).",
    "datalog_objects": "VI(R), VI_Request(R), VR_IT(R), VLink(R).",
    "datalog_relationships": "",
    "datalog_actions": "can_perform(S, R, Action) :- VIO(S), VI(R), Action = \"MLI:Request-VI\". can_perform(S, R, Action) :- VIO(S), VI(R), Action = \"MLI:Instantiate-VI\". can_perform(S, R, Action) :- VIO(S), VI(R), Action = \"MLI:Decommission-VI\". can_perform(S, R, Action) :- VIO(S), VI(R), Action = \"MLI:ReplanningVI:Modify-Time\". can_perform(S, R, Action) :- VIO(S), VI_Request(R), Action = \"MLI:Query-VI-Request-Status\". can_perform(S, R, Action) :- VIO(S), VI_Request(R), Action = \"MLI:Get-SLA-Offer\". can_perform(S, R, Action) :- VIO(S), VI_Request(R), Action = \"MLI:Sign-SLA-Offer\". can_perform(S, R, Action) :- VIO(S), VR_IT(R), Action = \"MLI:ReplanningVI:Add-VR-IT\". can_perform(S, R, Action) :- VIO(S), VR_IT(R), Action = \"MLI:ReplanningVI:Modify-VR-IT\". can_perform(S, R, Action) :- VIO(S), VR_IT(R), Action = \"MLI:ReplanningVI:Delete-VR-IT\". can_perform(S, R, Action) :- VIO(S), VLink(R), Action = \"MLI:ReplanningVI:Add-VLink\". can_perform(S, R, Action) :- VIO(S), VLink(R), Action = \"MLI:ReplanningVI:Modify-VLink\". can_perform(S, R, Action) :- VIO(S), VLink(R), Action = \"MLI:ReplanningVI:Delete-VLink\"."
}
'@
$ws.Range("C12").Value = $v_C12

$v_C13 = @'
{
    "datalog_subjects": "VIO(S).",
    "datalog_objects": "VR(R), VR_RP(RP), VR_RP_Info(RI), VR_State_Info(RSI), VR_Power_Info(RPI), VR_Status_Info(RSI2), VR_Info(RI2).",
    "datalog_relationships": "",
    "datalog_actions": "can_perform(S, R, A) :- VIO(S), VR(R), member(A, [\"SLI:Decommission-VR-IT\", \"SLI:Operate-VR\", \"SLI:Add-VirtualNetworkIf\", \"SLI:Remove-VirtualNetworkIf\", \"SLI:Create-StorageImage\", \"SLI:Remove-StorageImage\"]). can_perform(S, RP, \"SLI:Instantiate-VR-IT\") :- VIO(S), VR_RP(RP). can_perform(S, RI, \"SLI:Get-Available-VR-Pool-IT\") :- VIO(S), VR_RP_Info(RI). can_perform(S, RSI, \"SLI:Monitor-VR-Info\") :- VIO(S), VR_State_Info(RSI). can_perform(S, RPI, \"SLI:Monitor-VR-Info\") :- VIO(S), VR_Power_Info(RPI). can_perform(S, RSI2, \"SLI:Monitor-VR-Info\") :- VIO(S), VR_Status_Info(RSI2). can_perform(S, RI2, A) :- VIO(S), VR_Info(RI2), member(A, [\"SLI:Subscribe-VR-Monitoring\", \"SLI:Unsubscribe-VR-Monitoring\"])."
}
'@
$ws.Range("C13").Value = $v_C13

